$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 181317.5867854136
$ws.Cells.Item(3, 2).Value = 181529.4398680006
$ws.Cells.Item(4, 2).Value = 181744.2543182469
$ws.Cells.Item(5, 2).Value = 181962.0708740724
$ws.Cells.Item(6, 2).Value = 182182.9307946062
$ws.Cells.Item(7, 2).Value = 182406.8758644073
$ws.Cells.Item(8, 2).Value = 182633.9483975851
$ws.Cells.Item(9, 2).Value = 182864.1912418147
$ws.Cells.Item(10, 2).Value = 183097.6477822428
$ws.Cells.Item(11, 2).Value = 183334.3619452768
$ws.Cells.Item(12, 2).Value = 183574.3782022521
$ws.Cells.Item(13, 2).Value = 183817.7415729734
$ws.Cells.Item(14, 2).Value = 184064.4976291245
$ws.Cells.Item(15, 2).Value = 184314.6924975413
$ws.Cells.Item(16, 2).Value = 184568.3728633453
$ws.Cells.Item(17, 2).Value = 184825.5859729327
$ws.Cells.Item(18, 2).Value = 185086.3796368164
$ws.Cells.Item(19, 2).Value = 185350.8022323175
$ws.Cells.Item(20, 2).Value = 185618.9027061055
$ws.Cells.Item(21, 2).Value = 185890.7305765832
$ws.Cells.Item(22, 2).Value = 186166.3359361186
$ws.Cells.Item(23, 2).Value = 186445.7694531204
$ws.Cells.Item(24, 2).Value = 186729.0823739608
$ws.Cells.Item(25, 2).Value = 187016.3265247454
$ws.Cells.Item(26, 2).Value = 187307.5543129329
$ws.Cells.Item(27, 2).Value = 187602.8187288088
$ws.Cells.Item(28, 2).Value = 187902.1733468184
$ws.Cells.Item(29, 2).Value = 188205.6723267625
$ws.Cells.Item(30, 2).Value = 188513.3704148652
$ws.Cells.Item(31, 2).Value = 188825.3229447215
$ws.Cells.Item(32, 2).Value = 189141.5858381331
$ws.Cells.Item(33, 2).Value = 189462.2156058459
$ws.Cells.Item(34, 2).Value = 189787.2693481983
$ws.Cells.Item(35, 2).Value = 190116.8047556977
$ws.Cells.Item(36, 2).Value = 190450.8801095393
$ws.Cells.Item(37, 2).Value = 190789.5542820834
$ws.Cells.Item(38, 2).Value = 191132.8867373109
$ws.Cells.Item(39, 2).Value = 191480.9375312792
$ws.Cells.Item(40, 2).Value = 191833.7673125981
$ws.Cells.Item(41, 2).Value = 192191.4373229511
$ws.Cells.Item(42, 2).Value = 192554.009397688
$ws.Cells.Item(43, 2).Value = 192921.5459665173
$ws.Cells.Item(44, 2).Value = 193294.1100543253
$ws.Cells.Item(45, 2).Value = 193671.7652821573
$ws.Cells.Item(46, 2).Value = 194054.5758683912
$ws.Cells.Item(47, 2).Value = 194442.6066301388
$ws.Cells.Item(48, 2).Value = 194835.9229849133
$ws.Cells.Item(49, 2).Value = 195234.5909526021
$ws.Cells.Item(50, 2).Value = 195638.677157782
$ws.Cells.Item(51, 2).Value = 196048.2488324233
$ws.Cells.Item(52, 2).Value = 196463.3738190243
$ws.Cells.Item(53, 2).Value = 196884.1205742226
$ws.Cells.Item(54, 2).Value = 197310.55817293
$ws.Cells.Item(55, 2).Value = 197742.7563130389
$ws.Cells.Item(56, 2).Value = 198180.7853207525
$ws.Cells.Item(57, 2).Value = 198624.7161565858
$ws.Cells.Item(58, 2).Value = 199074.6204220951
$ws.Cells.Item(59, 2).Value = 199530.5703673827
$ws.Cells.Item(60, 2).Value = 199992.6388994374
$ws.Cells.Item(61, 2).Value = 200460.89959136
$ws.Cells.Item(62, 2).Value = 200935.4266925335
$ws.Cells.Item(63, 2).Value = 201416.2951397895
$ws.Cells.Item(64, 2).Value = 201903.5805696307
$ws.Cells.Item(65, 2).Value = 202397.3593315603
$ws.Cells.Item(66, 2).Value = 202897.7085025804
$ws.Cells.Item(67, 2).Value = 203404.7059029059
$ws.Cells.Item(68, 2).Value = 203918.4301129563
$ws.Cells.Item(69, 2).Value = 204438.960491672
$ws.Cells.Item(70, 2).Value = 204966.3771962127
$ws.Cells.Item(71, 2).Value = 205500.7612030854
$ws.Cells.Item(72, 2).Value = 206042.19433075
$ws.Cells.Item(73, 2).Value = 206590.7592637533
$ws.Cells.Item(74, 2).Value = 207146.5395784328
$ws.Cells.Item(75, 2).Value = 207709.6197702347
$ws.Cells.Item(76, 2).Value = 208280.0852826841
$ws.Cells.Item(77, 2).Value = 208858.0225380453
$ws.Cells.Item(78, 2).Value = 209443.5189697051
$ws.Cells.Item(79, 2).Value = 210036.6630563076
$ws.Cells.Item(80, 2).Value = 210637.5443576697
$ws.Cells.Item(81, 2).Value = 211246.2535524949
$ws.Cells.Item(82, 2).Value = 211862.8824779056
$ws.Cells.Item(83, 2).Value = 212487.5241708073
$ws.Cells.Item(84, 2).Value = 213120.2729110909
$ws.Cells.Item(85, 2).Value = 213761.2242666786
$ws.Cells.Item(86, 2).Value = 214410.4751404097
$ws.Cells.Item(87, 2).Value = 215068.1238187588
$ws.Cells.Item(88, 2).Value = 215734.2700223756
$ws.Cells.Item(89, 2).Value = 216409.0149584261
$ws.Cells.Item(90, 2).Value = 217092.4613747112
$ws.Cells.Item(91, 2).Value = 217784.7136155283
$ws.Cells.Item(92, 2).Value = 218485.8776792498
$ws.Cells.Item(93, 2).Value = 219196.0612775628
$ws.Cells.Item(94, 2).Value = 219915.3738963291
$ws.Cells.Item(95, 2).Value = 220643.9268580101
$ws.Cells.Item(96, 2).Value = 221381.8333855907
$ws.Cells.Item(97, 2).Value = 222129.2086679398
$ws.Cells.Item(98, 2).Value = 222886.1699265296
$ws.Cells.Item(99, 2).Value = 223652.8364834361
$ws.Cells.Item(100, 2).Value = 224429.3298305372
$ws.Cells.Item(101, 2).Value = 225215.7736998086
$ws.Cells.Item(102, 2).Value = 226012.2941346318
$ws.Cells.Item(103, 2).Value = 226819.0195620042
$ws.Cells.Item(104, 2).Value = 227636.0808655439
$ws.Cells.Item(105, 2).Value = 228463.6114591803
$ws.Cells.Item(106, 2).Value = 229301.7473614094
$ws.Cells.Item(107, 2).Value = 230150.627269998
$ws.Cells.Item(108, 2).Value = 231010.3926370013
$ws.Cells.Item(109, 2).Value = 231881.1877439782
$ws.Cells.Item(110, 2).Value = 232763.1597772602
$ws.Cells.Item(111, 2).Value = 233656.4589031481
$ws.Cells.Item(112, 2).Value = 234561.2383428939
$ws.Cells.Item(113, 2).Value = 235477.6544473321
$ws.Cells.Item(114, 2).Value = 236405.8667710198
$ws.Cells.Item(115, 2).Value = 237346.0381457449
$ws.Cells.Item(116, 2).Value = 238298.3347532626
$ws.Cells.Item(117, 2).Value = 239262.9261971179
$ws.Cells.Item(118, 2).Value = 240239.9855734162
$ws.Cells.Item(119, 2).Value = 241229.6895404062
$ws.Cells.Item(120, 2).Value = 242232.2183867345
$ws.Cells.Item(121, 2).Value = 243247.7560982422
$ws.Cells.Item(122, 2).Value = 244276.4904231785
$ws.Cells.Item(123, 2).Value = 245318.6129356939
$ws.Cells.Item(124, 2).Value = 246374.3190975049
$ws.Cells.Item(125, 2).Value = 247443.8083176071
$ws.Cells.Item(126, 2).Value = 248527.2840099283
$ws.Cells.Item(127, 2).Value = 249624.953648819
$ws.Cells.Item(128, 2).Value = 250737.0288222828
$ws.Cells.Item(129, 2).Value = 251863.7252828582
$ws.Cells.Item(130, 2).Value = 253005.2629960703
$ws.Cells.Item(131, 2).Value = 254161.8661863782
$ws.Cells.Item(132, 2).Value = 255333.7633805575
$ws.Cells.Item(133, 2).Value = 256521.187448464
$ws.Cells.Item(134, 2).Value = 257724.3756411286
$ws.Cells.Item(135, 2).Value = 258943.5696261571
$ws.Cells.Item(136, 2).Value = 260179.0155204086
$ws.Cells.Item(137, 2).Value = 261430.9639199371
$ws.Cells.Item(138, 2).Value = 262699.6699272027
$ws.Cells.Item(139, 2).Value = 263985.393175557
$ws.Cells.Item(140, 2).Value = 265288.397851028
$ws.Cells.Item(141, 2).Value = 266608.9527114388
$ws.Cells.Item(142, 2).Value = 267947.3311029088
$ws.Cells.Item(143, 2).Value = 269303.8109737896
$ws.Cells.Item(144, 2).Value = 270678.6748861194
$ws.Cells.Item(145, 2).Value = 272072.2100246648
$ws.Cells.Item(146, 2).Value = 273484.7082036565
$ws.Cells.Item(147, 2).Value = 274916.4658713228
$ws.Cells.Item(148, 2).Value = 276367.7841123415
$ws.Cells.Item(149, 2).Value = 277838.9686483385
$ws.Cells.Item(150, 2).Value = 279330.329836586
$ws.Cells.Item(151, 2).Value = 280842.1826670404
$ws.Cells.Item(152, 2).Value = 282374.8467579001
$ws.Cells.Item(153, 2).Value = 283928.6463498495
$ws.Cells.Item(154, 2).Value = 285503.9102991798
$ws.Cells.Item(155, 2).Value = 287100.9720699842
$ws.Cells.Item(156, 2).Value = 288720.1697256276
$ws.Cells.Item(157, 2).Value = 290361.8459197137
$ws.Cells.Item(158, 2).Value = 292026.347886766
$ws.Cells.Item(159, 2).Value = 293714.0274328554
$ws.Cells.Item(160, 2).Value = 295425.2409264144
$ws.Cells.Item(161, 2).Value = 297160.3492894782
$ws.Cells.Item(162, 2).Value = 298919.7179896108
$ws.Cells.Item(163, 2).Value = 300703.7170327631
$ws.Cells.Item(164, 2).Value = 302512.7209573313
$ws.Cells.Item(165, 2).Value = 304347.108829676
$ws.Cells.Item(166, 2).Value = 306207.2642413719
$ws.Cells.Item(167, 2).Value = 308093.5753084606
$ws.Cells.Item(168, 2).Value = 310006.4346729766
$ws.Cells.Item(169, 2).Value = 311946.2395070214
$ws.Cells.Item(170, 2).Value = 313913.3915196662
$ws.Cells.Item(171, 2).Value = 315908.2969669477
$ws.Cells.Item(172, 2).Value = 317931.3666652448
$ws.Cells.Item(173, 2).Value = 319983.0160082993
$ws.Cells.Item(174, 2).Value = 322063.6649881576
$ws.Cells.Item(175, 2).Value = 324173.7382203019
$ws.Cells.Item(176, 2).Value = 326313.6649732371
$ws.Cells.Item(177, 2).Value = 328483.879202798
$ws.Cells.Item(178, 2).Value = 330684.8195914319
$ws.Cells.Item(179, 2).Value = 332916.9295927187
$ws.Cells.Item(180, 2).Value = 335180.6574813704
$ws.Cells.Item(181, 2).Value = 337476.4564089612
$ws.Cells.Item(182, 2).Value = 339804.7844656258
$ws.Cells.Item(183, 2).Value = 342166.1047479551
$ws.Cells.Item(184, 2).Value = 344560.8854333233
$ws.Cells.Item(185, 2).Value = 346989.5998608633
$ws.Cells.Item(186, 2).Value = 349452.7266193044
$ws.Cells.Item(187, 2).Value = 351950.7496418806
$ws.Cells.Item(188, 2).Value = 354484.1583085125
$ws.Cells.Item(189, 2).Value = 357053.44755545
$ws.Cells.Item(190, 2).Value = 359659.11799257
$ws.Cells.Item(191, 2).Value = 362301.6760285049
$ws.Cells.Item(192, 2).Value = 364981.6340037716
$ws.Cells.Item(193, 2).Value = 367699.5103320757
$ws.Cells.Item(194, 2).Value = 370455.8296499404
$ws.Cells.Item(195, 2).Value = 373251.1229748193
$ws.Cells.Item(196, 2).Value = 376085.927871837
$ws.Cells.Item(197, 2).Value = 378960.7886292958
$ws.Cells.Item(198, 2).Value = 381876.2564430892
$ws.Cells.Item(199, 2).Value = 384832.8896101413
$ws.Cells.Item(200, 2).Value = 387831.2537310043
$ws.Cells.Item(201, 2).Value = 390871.9219217273
$ws.Cells.Item(202, 2).Value = 393955.4750351128
$ws.Cells.Item(203, 2).Value = 397082.501891468
$ws.Cells.Item(204, 2).Value = 400253.5995189323
$ws.Cells.Item(205, 2).Value = 403469.3734035591
$ws.Cells.Item(206, 2).Value = 406730.4377491747
$ws.Cells.Item(207, 2).Value = 410037.4157471403
$ws.Cells.Item(208, 2).Value = 413390.9398561041
$ws.Cells.Item(209, 2).Value = 416791.6520918431
$ws.Cells.Item(210, 2).Value = 420240.204327276
$ws.Cells.Item(211, 2).Value = 423737.2586027455
$ws.Cells.Item(212, 2).Value = 427283.487446653
$ws.Cells.Item(213, 2).Value = 430879.5742065363
$ws.Cells.Item(214, 2).Value = 434526.2133906751
$ws.Cells.Item(215, 2).Value = 438224.1110203149
$ws.Cells.Item(216, 2).Value = 441973.9849925995
$ws.Cells.Item(217, 2).Value = 445776.5654542907
$ws.Cells.Item(218, 2).Value = 449632.595186382
$ws.Cells.Item(219, 2).Value = 453542.8299996743
$ws.Cells.Item(220, 2).Value = 457508.0391414246
$ws.Cells.Item(221, 2).Value = 461529.0057131489
$ws.Cells.Item(222, 2).Value = 465606.5270996726
$ws.Cells.Item(223, 2).Value = 469741.4154095265
$ws.Cells.Item(224, 2).Value = 473934.4979267766
$ws.Cells.Item(225, 2).Value = 478186.61757439
$ws.Cells.Item(226, 2).Value = 482498.6333892205
$ws.Cells.Item(227, 2).Value = 486871.4210087234
$ws.Cells.Item(228, 2).Value = 491305.8731694774
$ws.Cells.Item(229, 2).Value = 495802.9002176187
$ws.Cells.Item(230, 2).Value = 500363.4306312723
$ws.Cells.Item(231, 2).Value = 504988.4115550739
$ws.Cells.Item(232, 2).Value = 509678.8093468602
$ws.Cells.Item(233, 2).Value = 514435.6101366181
$ws.Cells.Item(234, 2).Value = 519259.820397765
$ws.Cells.Item(235, 2).Value = 524152.4675308322
$ws.Cells.Item(236, 2).Value = 529114.6004596144
$ws.Cells.Item(237, 2).Value = 534147.2902398458
$ws.Cells.Item(238, 2).Value = 539251.6306804394
$ws.Cells.Item(239, 2).Value = 544428.7389773368
$ws.Cells.Item(240, 2).Value = 549679.7563599845
$ws.Cells.Item(241, 2).Value = 555005.8487504424
$ws.Cells.Item(242, 2).Value = 560408.2074351287
$ws.Cells.Item(243, 2).Value = 565888.0497491665
$ws.Cells.Item(244, 2).Value = 571446.6197732913
$ws.Cells.Item(245, 2).Value = 577085.1890432568
$ws.Cells.Item(246, 2).Value = 582805.0572716413
$ws.Cells.Item(247, 2).Value = 588607.5530819532
$ws.Cells.Item(248, 2).Value = 594494.034754876
$ws.Cells.Item(249, 2).Value = 600465.8909864889
$ws.Cells.Item(250, 2).Value = 606524.5416582463
$ws.Cells.Item(251, 2).Value = 612671.4386184823
$ws.Cells.Item(252, 2).Value = 618908.0664751426
$ws.Cells.Item(253, 2).Value = 625235.9433994255
$ws.Cells.Item(254, 2).Value = 631656.6219399656
$ws.Cells.Item(255, 2).Value = 638171.6898471242
$ws.Cells.Item(256, 2).Value = 644782.7709069295
$ws.Cells.Item(257, 2).Value = 651491.5257841346
$ws.Cells.Item(258, 2).Value = 658299.6528738079
$ws.Cells.Item(259, 2).Value = 665208.8891608036
$ws.Cells.Item(260, 2).Value = 672221.0110864117
$ws.Cells.Item(261, 2).Value = 679337.835421382
$ws.Cells.Item(262, 2).Value = 686561.2201444913
$ws.Cells.Item(263, 2).Value = 693893.0653256988
$ws.Cells.Item(264, 2).Value = 701335.314012876
$ws.Cells.Item(265, 2).Value = 708889.9531210224
$ws.Cells.Item(266, 2).Value = 716559.0143227541
$ws.Cells.Item(267, 2).Value = 724344.5749387913
$ws.Cells.Item(268, 2).Value = 732248.758827064
$ws.Cells.Item(269, 2).Value = 740273.7372689296
$ws.Cells.Item(270, 2).Value = 748421.7298509241
$ws.Cells.Item(271, 2).Value = 756695.0053403305
$ws.Cells.Item(272, 2).Value = 765095.8825527519
$ws.Cells.Item(273, 2).Value = 773626.7312097271
$ws.Cells.Item(274, 2).Value = 782289.9727843517
$ws.Cells.Item(275, 2).Value = 791088.0813326794
$ws.Cells.Item(276, 2).Value = 800023.5843085918
$ws.Cells.Item(277, 2).Value = 809099.0633596792
$ws.Cells.Item(278, 2).Value = 818317.1551015048
$ws.Cells.Item(279, 2).Value = 827680.5518675356
$ws.Cells.Item(280, 2).Value = 837192.0024318154
$ws.Cells.Item(281, 2).Value = 846854.312701356
$ws.Cells.Item(282, 2).Value = 856670.346375039
$ws.Cells.Item(283, 2).Value = 866643.0255656763
$ws.Cells.Item(284, 2).Value = 876775.3313817157
$ws.Cells.Item(285, 2).Value = 887070.3044649207
$ws.Cells.Item(286, 2).Value = 897531.0454801938
$ws.Cells.Item(287, 2).Value = 908160.7155535422
$ws.Cells.Item(288, 2).Value = 918962.536654025
$ws.Cells.Item(289, 2).Value = 929939.791915379
$ws.Cells.Item(290, 2).Value = 941095.8258928105
$ws.Cells.Item(291, 2).Value = 952434.0447503353
$ws.Cells.Item(292, 2).Value = 963957.9163738508
$ws.Cells.Item(293, 2).Value = 975670.9704049787
$ws.Cells.Item(294, 2).Value = 987576.7981905874
$ws.Cells.Item(295, 2).Value = 999679.0526427121
$ws.Cells.Item(296, 2).Value = 1011981.448003486
$ws.Cells.Item(297, 2).Value = 1024487.759509552
$ws.Cells.Item(298, 2).Value = 1037201.822950287
$ws.Cells.Item(299, 2).Value = 1050127.534114054
$ws.Cells.Item(300, 2).Value = 1063268.848116598
$ws.Cells.Item(301, 2).Value = 1076629.778605609
$ws.Cells.Item(302, 2).Value = 1090214.396835345
$ws.Cells.Item(303, 2).Value = 1257115.036895679
$ws.Cells.Item(304, 2).Value = 1426610.529747827
$ws.Cells.Item(305, 2).Value = 1598744.916083281
$ws.Cells.Item(306, 2).Value = 1773562.785463792
$ws.Cells.Item(307, 2).Value = 1951109.261696993
$ws.Cells.Item(308, 2).Value = 2131429.986443242
$ws.Cells.Item(309, 2).Value = 2314571.100985378
$ws.Cells.Item(310, 2).Value = 2500579.226093488
$ws.Cells.Item(311, 2).Value = 2689501.439917476
$ws.Cells.Item(312, 2).Value = 2881385.253841698
$ws.Cells.Item(313, 2).Value = 3076278.586237024
$ws.Cells.Item(314, 2).Value = 3274229.734047352
$ws.Cells.Item(315, 2).Value = 3475287.342149931
$ws.Cells.Item(316, 2).Value = 3679500.370431011
$ws.Cells.Item(317, 2).Value = 3886918.058521098
$ws.Cells.Item(318, 2).Value = 4097589.888137488
$ws.Cells.Item(319, 2).Value = 4311565.542984924
$ws.Cells.Item(320, 2).Value = 4528894.866169413
$ws.Cells.Item(321, 2).Value = 4749627.815084763
$ws.Cells.Item(322, 2).Value = 4973814.413735781
$ws.Cells.Item(323, 2).Value = 5201504.702467853
$ws.Cells.Item(324, 2).Value = 5432748.685077859
$ws.Cells.Item(325, 2).Value = 5667596.273288168
$ws.Cells.Item(326, 2).Value = 5906097.22857149
$ws.Cells.Item(327, 2).Value = 6148301.101321952
$ws.Cells.Item(328, 2).Value = 6394257.167375253
$ws.Cells.Item(329, 2).Value = 6644014.361888752
$ws.Cells.Item(330, 2).Value = 6897621.210600981
$ws.Cells.Item(331, 2).Value = 7155125.758499053
$ws.Cells.Item(332, 2).Value = 7416575.495931973
$ws.Cells.Item(333, 2).Value = 7682017.282217473
$ws.Cells.Item(334, 2).Value = 7951497.266800604
$ws.Cells.Item(335, 2).Value = 8225060.808032771
$ws.Cells.Item(336, 2).Value = 8502752.389650924
$ws.Cells.Item(337, 2).Value = 8784615.535048174
$ws.Cells.Item(338, 2).Value = 9070692.719438449
$ws.Cells.Item(339, 2).Value = 9361025.280029736
$ws.Cells.Item(340, 2).Value = 9655653.324332999
$ws.Cells.Item(341, 2).Value = 9954615.636745347
$ws.Cells.Item(342, 2).Value = 10257949.58355894
$ws.Cells.Item(343, 2).Value = 10565691.01655947
$ws.Cells.Item(344, 2).Value = 10877874.17539001
$ws.Cells.Item(345, 2).Value = 11194531.58886902
$ws.Cells.Item(346, 2).Value = 11515693.97546258
$ws.Cells.Item(347, 2).Value = 11841390.14312382
$ws.Cells.Item(348, 2).Value = 12171646.88872316
$ws.Cells.Item(349, 2).Value = 12506488.89730519
$ws.Cells.Item(350, 2).Value = 12845938.64141782
$ws.Cells.Item(351, 2).Value = 13190016.2807707
$ws.Cells.Item(352, 2).Value = 13538739.56248846
$ws.Cells.Item(353, 2).Value = 13892123.7222342
$ws.Cells.Item(354, 2).Value = 14250181.38648617
$ws.Cells.Item(355, 2).Value = 14612922.47625818
$ws.Cells.Item(356, 2).Value = 14980354.11256073
$ws.Cells.Item(357, 2).Value = 15352480.52390479
$ws.Cells.Item(358, 2).Value = 15343870.05609418
$ws.Cells.Item(359, 2).Value = 15335007.29013677
$ws.Cells.Item(360, 2).Value = 15325879.6117861
$ws.Cells.Item(361, 2).Value = 15316473.78048362
$ws.Cells.Item(362, 2).Value = 15306775.90543684
$ws.Cells.Item(363, 2).Value = 15296771.42132337
$ws.Cells.Item(364, 2).Value = 15286445.06366115
$ws.Cells.Item(365, 2).Value = 15275780.84389077
$ws.Cells.Item(366, 2).Value = 15264762.02422051
$ws.Cells.Item(367, 2).Value = 15253371.09229142
$ws.Cells.Item(368, 2).Value = 15241589.73572542
$ws.Cells.Item(369, 2).Value = 15229398.81662627
$ws.Cells.Item(370, 2).Value = 15216778.34611044
$ws.Cells.Item(371, 2).Value = 15203707.45895199
$ws.Cells.Item(372, 2).Value = 15190164.38843382
$ws.Cells.Item(373, 2).Value = 15176126.4415057
$ws.Cells.Item(374, 2).Value = 15161569.97435804
$ws.Cells.Item(375, 2).Value = 15146470.36852993
$ws.Cells.Item(376, 2).Value = 15130802.00767867
$ws.Cells.Item(377, 2).Value = 15114538.25514836
$ws.Cells.Item(378, 2).Value = 15097651.43248526
$ws.Cells.Item(379, 2).Value = 15080112.7990575
$ws.Cells.Item(380, 2).Value = 15061892.53294807
$ws.Cells.Item(381, 2).Value = 15042959.71330057
$ws.Cells.Item(382, 2).Value = 15023282.30430817
$ws.Cells.Item(383, 2).Value = 15002827.14104753
$ws.Cells.Item(384, 2).Value = 14981559.91737066
$ws.Cells.Item(385, 2).Value = 14959445.17607833
$ws.Cells.Item(386, 2).Value = 14936446.30161009
$ws.Cells.Item(387, 2).Value = 14912525.51549609
$ws.Cells.Item(388, 2).Value = 14887643.87482646
$ws.Cells.Item(389, 2).Value = 14861761.27400326
$ws.Cells.Item(390, 2).Value = 14834836.45004973
$ws.Cells.Item(391, 2).Value = 14806826.99175879
$ws.Cells.Item(392, 2).Value = 14777689.35297093
$ws.Cells.Item(393, 2).Value = 14747378.87027708
$ws.Cells.Item(394, 2).Value = 14715849.78544696
$ws.Cells.Item(395, 2).Value = 14683055.27288686
$ws.Cells.Item(396, 2).Value = 14648947.47243186
$ws.Cells.Item(397, 2).Value = 14613477.5277778
$ws.Cells.Item(398, 2).Value = 14576595.63085521
$ws.Cells.Item(399, 2).Value = 14538251.07244312
$ws.Cells.Item(400, 2).Value = 14498392.29931314
$ws.Cells.Item(401, 2).Value = 14456966.9781845
$ws.Cells.Item(402, 2).Value = 14413922.06675763
$ws.Cells.Item(403, 2).Value = 14369203.89207827
$ws.Cells.Item(404, 2).Value = 14322758.2364648
$ws.Cells.Item(405, 2).Value = 14274530.43120908
$ws.Cells.Item(406, 2).Value = 14224465.45823498
$ws.Cells.Item(407, 2).Value = 14172508.05986926
$ws.Cells.Item(408, 2).Value = 14118602.85684608
$ws.Cells.Item(409, 2).Value = 14062694.47462919
$ws.Cells.Item(410, 2).Value = 14004727.67809511
$ws.Cells.Item(411, 2).Value = 13944647.51457571
$ws.Cells.Item(412, 2).Value = 13882399.46521039
$ws.Cells.Item(413, 2).Value = 13817929.60450576
$ws.Cells.Item(414, 2).Value = 13751184.76794523
$ws.Cells.Item(415, 2).Value = 13682112.72743219
$ws.Cells.Item(416, 2).Value = 13610662.37428804
$ws.Cells.Item(417, 2).Value = 13536783.90946218
$ws.Cells.Item(418, 2).Value = 13460429.04054305
$ws.Cells.Item(419, 2).Value = 13381551.1850909
$ws.Cells.Item(420, 2).Value = 13300105.67974112
$ws.Cells.Item(421, 2).Value = 13216049.99445559
$ws.Cells.Item(422, 2).Value = 13129343.95122634
$ws.Cells.Item(423, 2).Value = 13039949.94646353
$ws.Cells.Item(424, 2).Value = 12947833.17622731
$ws.Cells.Item(425, 2).Value = 12852961.86339264
$ws.Cells.Item(426, 2).Value = 12755307.48576729
$ws.Cells.Item(427, 2).Value = 12654845.00411722
$ws.Cells.Item(428, 2).Value = 12551553.08899103
$ws.Cells.Item(429, 2).Value = 12445414.34517705
$ws.Cells.Item(430, 2).Value = 12336415.53257319
$ws.Cells.Item(431, 2).Value = 12224547.78220279
$ws.Cells.Item(432, 2).Value = 12109806.80606864
$ws.Cells.Item(433, 2).Value = 11992193.09950405
$ws.Cells.Item(434, 2).Value = 11871712.13465497
$ws.Cells.Item(435, 2).Value = 11748374.5437102
$ws.Cells.Item(436, 2).Value = 11622196.2904902
$ws.Cells.Item(437, 2).Value = 11493198.82900771
$ws.Cells.Item(438, 2).Value = 11361409.24762686
$ws.Cells.Item(439, 2).Value = 11226860.39747216
$ws.Cells.Item(440, 2).Value = 11089591.00377422
$ws.Cells.Item(441, 2).Value = 10949645.75888623
$ws.Cells.Item(442, 2).Value = 10807075.39576454
$ws.Cells.Item(443, 2).Value = 10661936.74077653
$ws.Cells.Item(444, 2).Value = 10514292.74478212
$ws.Cells.Item(445, 2).Value = 10364212.49152775
$ws.Cells.Item(446, 2).Value = 10211771.1824972
$ws.Cells.Item(447, 2).Value = 10057050.09747834
$ws.Cells.Item(448, 2).Value = 9900136.530230129
$ws.Cells.Item(449, 2).Value = 9741123.698768903
$ws.Cells.Item(450, 2).Value = 9580110.629935747
$ws.Cells.Item(451, 2).Value = 9417202.018057304
$ws.Cells.Item(452, 2).Value = 9252508.057669174
